$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (QUILT): Progress 38 -> 50
$ws.Range("B4").Value = 50

# Row 9 (ZENITH (ALN-AGT01-008)): Progress 75 -> 100
$ws.Range("B9").Value = 100

$wb.Save()
